$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 2) below the header row, mirroring the new
# "Deposit / Crypto / ETH / 0" transaction entry.
$ws.Range("E2").Value = "Deposit"
$ws.Range("N2").Value = "Crypto"
$ws.Range("P2").Value = "ETH"
$ws.Range("T2").Value = 0

# Update the view: scroll so column O is the left-most visible column,
# and move the active selection to T3 (a single cell, no more range
# selected across the old data rows).
$win = $wb.Windows.Item(1)
$win.ScrollRow = 1
$win.ScrollColumn = 15
$ws.Range("T3").Select()
